$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")
$ws.Rows.Item(211).Insert()
$ws.Cells.Item(211, 1).Value = "trans"
$ws.Cells.Item(211, 2).Value = "FpUCD"
$ws.Cells.Item(211, 3).Value = "Fares per Unit Cargo Distance"
$ws.Cells.Item(211, 6).Value = "medium"
$ws.Cells.Item(218, 6).Copy()
$ws.Cells.Item(211, 6).PasteSpecial(-4122)

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()

Write-Host "done"
